$wb = $excel.ActiveWorkbook

# Rename sheets (tab names) to new task-order IDs
$wb.Worksheets.Item(1).Name = "GNG_TO-16502911679634035"
$wb.Worksheets.Item(2).Name = "NB_TO-16502911714782865"
$wb.Worksheets.Item(3).Name = "RS_TO-16502911714802582"
$wb.Worksheets.Item(4).Name = "TOL_TO-16502911715412538"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16502911716352825"

# Sheet 1 (GNG) - update stim filenames
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-1650291167933398.csv"
$ws1.Range("B3").Value = "GNG_stims-1650291167945397.csv"
$ws1.Range("B4").Value = "go_stims-1650291167946397.csv"
$ws1.Range("B5").Value = "GNG_stims-16502911679614322.csv"

# Sheet 2 (NB) - update stim filenames
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "ZB-match_1-16502911681981401.csv"
$ws2.Range("B3").Value = "OB-16502911692421424.csv"
$ws2.Range("B4").Value = "TB-16502911698031068.csv"
$ws2.Range("B5").Value = "TB-16502911714662495.csv"
$ws2.Range("B6").Value = "OB-16502911685731049.csv"
$ws2.Range("B7").Value = "ZB-match_0-16502911680303955.csv"
$ws2.Range("B8").Value = "ZB-match_6-16502911681101093.csv"
$ws2.Range("B9").Value = "OB-16502911693641357.csv"
$ws2.Range("B10").Value = "TB-1650291169514105.csv"

# Sheet 4 (TOL) - update stim filenames
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16502911714932532.csv"
$ws4.Range("B3").Value = "ZM_stims-16502911714812539.csv"
$ws4.Range("B4").Value = "MM_stims-1650291171524285.csv"
$ws4.Range("B5").Value = "ZM_stims-1650291171494251.csv"
$ws4.Range("B6").Value = "MM_stims-16502911715402524.csv"
$ws4.Range("B7").Value = "ZM_stims-16502911715252557.csv"

# Sheet 5 (vSAT) - update stim filenames
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16502911715712845.csv"
$ws5.Range("B3").Value = "vSAT_stims-16502911716192508.csv"
$ws5.Range("B4").Value = "vSAT_stims-16502911715892506.csv"
$ws5.Range("B5").Value = "SAT_stims-1650291171545252.csv"
